# Rebuild the "Sheet1" table with the new Phase/Name/Description/Units/Note/
# Distribution Type/Parameter 1-7/Lower Limit/Upper Limit/Step layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - remove the old A:J table entirely before laying out the new one.
$ws.Cells.Clear()

# ---- Header row (row 1) ----
$ws.Cells.Item(1, 1).Value = "Phase"
$ws.Cells.Item(1, 2).Value = "Name"
$ws.Cells.Item(1, 3).Value = "Description"
$ws.Cells.Item(1, 4).Value = "Units"
$ws.Cells.Item(1, 5).Value = "Note"
$ws.Cells.Item(1, 6).Value = "Distribution Type"
$ws.Cells.Item(1, 7).Value = "Parameter 1"
$ws.Cells.Item(1, 8).Value = "Parameter 2"
$ws.Cells.Item(1, 9).Value = "Parameter 3"
$ws.Cells.Item(1, 10).Value = "Parameter 4"
$ws.Cells.Item(1, 11).Value = "Parameter 5"
$ws.Cells.Item(1, 12).Value = "Parameter 6"
$ws.Cells.Item(1, 13).Value = "Parameter 7"
$ws.Cells.Item(1, 14).Value = "Lower Limit"
$ws.Cells.Item(1, 15).Value = "Upper Limit"
$ws.Cells.Item(1, 16).Value = "Step"

# ---- Row 2 : Test Parameter 1 ----
$ws.Cells.Item(2, 1).Value = "Indoor;Underground;Outdoor"
$ws.Cells.Item(2, 2).Value = "Test Parameter 1"
$ws.Cells.Item(2, 3).Value = "The first test parameter"
$ws.Cells.Item(2, 4).Value = "Some Units"
$ws.Cells.Item(2, 6).Value = "Constant"
$ws.Cells.Item(2, 7).Value = 1
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 4
$ws.Cells.Item(2, 16).Value = 1

# ---- Row 3 : Test Parameter 2 ----
$ws.Cells.Item(3, 1).Value = "Indoor;Underground"
$ws.Cells.Item(3, 2).Value = "Test Parameter 2"
$ws.Cells.Item(3, 3).Value = "The second test parameter"
$ws.Cells.Item(3, 4).Value = "Some Units"
$ws.Cells.Item(3, 6).Value = "Constant"
$ws.Cells.Item(3, 7).Value = 2
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(3, 15).Value = 4
$ws.Cells.Item(3, 16).Value = 1

# ---- Row 4 : Test Parameter 3 ----
$ws.Cells.Item(4, 1).Value = "Outdoor"
$ws.Cells.Item(4, 2).Value = "Test Parameter 3"
$ws.Cells.Item(4, 3).Value = "The third test parameter"
$ws.Cells.Item(4, 4).Value = "Some Units"
$ws.Cells.Item(4, 6).Value = "Constant"
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 15).Value = 4
$ws.Cells.Item(4, 16).Value = 1

# Resize columns to fit the new content, like Excel would after editing.
$ws.Columns.AutoFit() | Out-Null

# Match the saved selection/active cell from the authored workbook.
$ws.Range("K5").Select() | Out-Null
